$wb = $excel.ActiveWorkbook

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7969.615
$ws.Range("I62").Value = 6633.3335
$ws.Range("J62").Value = 9115
$ws.Range("K62").Value = 6633.3335
$ws.Range("L62").Value = 9115
$ws.Range("M62").Value = -6009.3335
$ws.Range("N62").Value = -10363

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7969.615
$ws.Range("I65").Value = 6633.3335
$ws.Range("J65").Value = 9115
$ws.Range("K65").Value = 33166.6675
$ws.Range("L65").Value = 45575
$ws.Range("M65").Value = -30046.6675
$ws.Range("N65").Value = -51815

# Sheet ALC, row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2742.7812
$ws.Range("I76").Value = 2780.3914
$ws.Range("K76").Value = 2780.3914
$ws.Range("M76").Value = -2465.3914

# Sheet ALC, row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2742.7812
$ws.Range("I79").Value = 2780.3914
$ws.Range("K79").Value = 2780.3914
$ws.Range("M79").Value = -1688.3914

# Sheet ALC, row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 7600
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594

# Sheet ALC, row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 7600
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596

# Sheet ALC, row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 283.75
$ws.Range("I127").Value = 283.75
$ws.Range("K127").Value = 851.25
$ws.Range("M127").Value = 4108.75

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1110.375
$ws.Range("I129").Value = 513.1667
$ws.Range("J129").Value = 1248.1923
$ws.Range("K129").Value = 1539.5001
$ws.Range("L129").Value = 3744.5769
$ws.Range("M129").Value = 3460.4999
$ws.Range("N129").Value = -13744.5769

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7399.7075
$ws.Range("I132").Value = 6714.5557
$ws.Range("K132").Value = 20143.6671
$ws.Range("M132").Value = -17613.6671

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1387.7091
$ws.Range("I137").Value = 1586.9032
$ws.Range("J137").Value = 1130.4166
$ws.Range("K137").Value = 4760.7096
$ws.Range("L137").Value = 3391.2498
$ws.Range("M137").Value = -2210.7096
$ws.Range("N137").Value = -8491.2498

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 16177.637
$ws.Range("I141").Value = 6348.3335
$ws.Range("J141").Value = 27972.8
$ws.Range("K141").Value = 19045.0005
$ws.Range("L141").Value = 83918.39999999999
$ws.Range("M141").Value = -13865.0005
$ws.Range("N141").Value = -94278.39999999999

# Sheet ARM, row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 220.55556
$ws.Range("I4").Value = 227.14285
$ws.Range("J4").Value = 197.5
$ws.Range("K4").Value = 227.14285
$ws.Range("L4").Value = 197.5
$ws.Range("M4").Value = -111.14285
$ws.Range("N4").Value = -429.5

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3470.7932
$ws.Range("I61").Value = 3987.5715
$ws.Range("J61").Value = 2114.25
$ws.Range("K61").Value = 3987.5715
$ws.Range("L61").Value = 2114.25
$ws.Range("M61").Value = -3775.5715
$ws.Range("N61").Value = -2538.25

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2045.579
$ws.Range("I74").Value = 1554.875
$ws.Range("J74").Value = 4662.6665
$ws.Range("K74").Value = 1554.875
$ws.Range("L74").Value = 4662.6665
$ws.Range("M74").Value = -680.875
$ws.Range("N74").Value = -6410.6665

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2045.579
$ws.Range("I77").Value = 1554.875
$ws.Range("J77").Value = 4662.6665
$ws.Range("K77").Value = 7774.375
$ws.Range("L77").Value = 23313.3325
$ws.Range("M77").Value = -3406.375
$ws.Range("N77").Value = -32049.3325

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5085.5
$ws.Range("I132").Value = 2992.9473
$ws.Range("J132").Value = 6814.1304
$ws.Range("K132").Value = 8978.841899999999
$ws.Range("L132").Value = 20442.3912
$ws.Range("M132").Value = -6448.841899999999
$ws.Range("N132").Value = -25502.3912

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3470.7932
$ws.Range("I136").Value = 3987.5715
$ws.Range("J136").Value = 2114.25
$ws.Range("K136").Value = 11962.7145
$ws.Range("L136").Value = 6342.75
$ws.Range("M136").Value = -9412.7145
$ws.Range("N136").Value = -11442.75

# Sheet BSM, row 57
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 59633.332
$ws.Range("J57").Value = 59633.332
$ws.Range("L57").Value = 59633.332
$ws.Range("N57").Value = -61073.332

# Sheet BSM, row 136
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 59633.332
$ws.Range("J136").Value = 59633.332
$ws.Range("L136").Value = 59633.332
$ws.Range("N136").Value = -69833.33199999999

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10103391
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10103391
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10103391
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -10103981

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10103391
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10103391
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10103391
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -10103795

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4700.212
$ws.Range("I58").Value = 2040
$ws.Range("J58").Value = 7526.6875
$ws.Range("K58").Value = 2040
$ws.Range("L58").Value = 7526.6875
$ws.Range("M58").Value = -1837
$ws.Range("N58").Value = -7932.6875

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2846.7856
$ws.Range("I132").Value = 2201.5
$ws.Range("J132").Value = 3205.2778
$ws.Range("K132").Value = 6604.5
$ws.Range("L132").Value = 9615.8334
$ws.Range("M132").Value = -4074.5
$ws.Range("N132").Value = -14675.8334

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1579.356
$ws.Range("I134").Value = 1004
$ws.Range("J134").Value = 2790.6316
$ws.Range("K134").Value = 3012
$ws.Range("L134").Value = 8371.8948
$ws.Range("M134").Value = -477
$ws.Range("N134").Value = -13441.8948

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4700.212
$ws.Range("I136").Value = 2040
$ws.Range("J136").Value = 7526.6875
$ws.Range("K136").Value = 6120
$ws.Range("L136").Value = 22580.0625
$ws.Range("M136").Value = -3570
$ws.Range("N136").Value = -27680.0625

# Sheet CUL, row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1427.2858
$ws.Range("I98").Value = 777.7778
$ws.Range("J98").Value = 1914.4166
$ws.Range("K98").Value = 2333.3334
$ws.Range("L98").Value = 5743.2498
$ws.Range("M98").Value = -835.3334
$ws.Range("N98").Value = -8739.2498

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2739.1765
$ws.Range("I122").Value = 708
$ws.Range("K122").Value = 6372
$ws.Range("M122").Value = -3922

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1195.6487
$ws.Range("I131").Value = 1280.5454
$ws.Range("J131").Value = 1159.7307
$ws.Range("K131").Value = 3841.6362
$ws.Range("L131").Value = 3479.1921
$ws.Range("M131").Value = 1198.3638
$ws.Range("N131").Value = -13559.1921

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3145.889
$ws.Range("I80").Value = 3200
$ws.Range("J80").Value = 3139.125
$ws.Range("K80").Value = 3200
$ws.Range("L80").Value = 3139.125
$ws.Range("M80").Value = -2202
$ws.Range("N80").Value = -5135.125

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3145.889
$ws.Range("I83").Value = 3200
$ws.Range("J83").Value = 3139.125
$ws.Range("K83").Value = 16000
$ws.Range("L83").Value = 15695.625
$ws.Range("M83").Value = -11008
$ws.Range("N83").Value = -25679.625

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1242.258
$ws.Range("I102").Value = 975.4783
$ws.Range("J102").Value = 2009.25
$ws.Range("K102").Value = 975.4783
$ws.Range("L102").Value = 2009.25
$ws.Range("M102").Value = 646.5217
$ws.Range("N102").Value = -5253.25

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 713.5862
$ws.Range("I107").Value = 695.1111
$ws.Range("J107").Value = 743.8182
$ws.Range("K107").Value = 695.1111
$ws.Range("L107").Value = 743.8182
$ws.Range("M107").Value = 1224.8889
$ws.Range("N107").Value = -4583.8182

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2451.5312
$ws.Range("I132").Value = 1834.4375
$ws.Range("J132").Value = 3068.625
$ws.Range("K132").Value = 5503.3125
$ws.Range("L132").Value = 9205.875
$ws.Range("M132").Value = -2973.3125
$ws.Range("N132").Value = -14265.875

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2113.5625
$ws.Range("I122").Value = 760.8
$ws.Range("J122").Value = 4368.1665
$ws.Range("K122").Value = 2282.4
$ws.Range("L122").Value = 13104.4995
$ws.Range("M122").Value = 167.6000000000004
$ws.Range("N122").Value = -18004.4995

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1629.2812
$ws.Range("I132").Value = 1224.619
$ws.Range("J132").Value = 2401.818
$ws.Range("K132").Value = 3673.857
$ws.Range("L132").Value = 7205.454000000001
$ws.Range("M132").Value = -1143.857
$ws.Range("N132").Value = -12265.454
